# Applies the Mon Jul 17 18:32:39 UTC 2023 "Updated cryptos list" refresh:
# new Price (D) / Volume(1h) (E) quotes for every coin row, plus the
# WrappedliquidstakedEther2.0 <-> Dai row swap (rows 21 & 22).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold plain numeric-looking text (e.g. "1.001",
# "30.063.78") in the source workbook - t="inlineStr", no number format.
# Excel auto-converts a bare numeric string assigned to .Value into a real
# Number, which would change the stored cell type. Prefixing the value with
# a literal apostrophe (the normal Excel "force text" input) keeps it text,
# exactly like the source cell - this is done for every D-column value that
# would otherwise parse as a number; multi-dot strings like "2.171.68" are
# already unambiguous text and do not need it.
function Set-TextValue($range, [string]$value) {
    $ws.Range($range).Value = "'" + $value
}

$ws.Range("D2").Value = '30.063.78'
$ws.Range("E2").Value = '  -0.94%  '

$ws.Range("D3").Value = '1.898.84'
$ws.Range("E3").Value = '  -1.92%  '

Set-TextValue "D4" '1.001'
$ws.Range("E4").Value = '  -0.04%  '

Set-TextValue "D5" '0.7288'
$ws.Range("E5").Value = '  -6.01%  '

Set-TextValue "D6" '242.59'
$ws.Range("E6").Value = '  -1.49%  '

Set-TextValue "D7" '1.002'
$ws.Range("E7").Value = '  +0.15%  '

Set-TextValue "D8" '0.3096'
$ws.Range("E8").Value = '  -3.47%  '

$ws.Range("E9").Value = '  -6.11%  '

Set-TextValue "D10" '0.06896'
$ws.Range("E10").Value = '  -2.38%  '

Set-TextValue "D11" '0.7698'
$ws.Range("E11").Value = '  -1.52%  '

Set-TextValue "D12" '0.07943'
$ws.Range("E12").Value = '  -0.85%  '

$ws.Range("D13").Value = '1.897.31'
$ws.Range("E13").Value = '  -1.82%  '

Set-TextValue "D14" '5.253'
$ws.Range("E14").Value = '  -2.08%  '

Set-TextValue "D15" '90.90'
$ws.Range("E15").Value = '  -4.11%  '

$ws.Range("D16").Value = '30.045.88'
$ws.Range("E16").Value = '  -0.95%  '

Set-TextValue "D17" '14.14'
$ws.Range("E17").Value = '  -2.72%  '

Set-TextValue "D18" '5.748'
$ws.Range("E18").Value = '  -1.30%  '

Set-TextValue "D19" '0.000007757'
$ws.Range("E19").Value = '  -2.81%  '

Set-TextValue "D20" '237.08'
$ws.Range("E20").Value = '  -7.17%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.171.68'
$ws.Range("E21").Value = '  -0.57%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D22" '1.001'
$ws.Range("E22").Value = '  +0.27%  '

Set-TextValue "D23" '1.001'
$ws.Range("E23").Value = '  -0.06%  '

Set-TextValue "D24" '6.891'
$ws.Range("E24").Value = '  +1.97%  '

Set-TextValue "D25" '9.311'
$ws.Range("E25").Value = '  -2.72%  '

Set-TextValue "D26" '165.73'
$ws.Range("E26").Value = '  +1.33%  '

Set-TextValue "D27" '18.89'
$ws.Range("E27").Value = '  -1.09%  '

Set-TextValue "D28" '0.1274'
$ws.Range("E28").Value = '  -5.98%  '

Set-TextValue "D29" '2.019'
$ws.Range("E29").Value = '  -11.54%  '

$ws.Range("E30").Value = '  -1.47%  '

Set-TextValue "D31" '1.536'
$ws.Range("E31").Value = '  +1.15%  '

Set-TextValue "D32" '4.295'
$ws.Range("E32").Value = '  -3.01%  '

Set-TextValue "D33" '4.069'
$ws.Range("E33").Value = '  -1.53%  '

Set-TextValue "D34" '0.05083'
$ws.Range("E34").Value = '  -1.76%  '

Set-TextValue "D35" '1.270'
$ws.Range("E35").Value = '  -1.11%  '

Set-TextValue "D36" '0.7348'
$ws.Range("E36").Value = '  -2.16%  '

Set-TextValue "D37" '2.731'
$ws.Range("E37").Value = '  -1.48%  '

Set-TextValue "D38" '0.01918'
$ws.Range("E38").Value = '  -2.24%  '

$ws.Range("E39").Value = '  -1.26%  '

Set-TextValue "D40" '6.337'
$ws.Range("E40").Value = '  -1.69%  '

Set-TextValue "D41" '74.64'
$ws.Range("E41").Value = '  -5.30%  '

Set-TextValue "D42" '0.4430'
$ws.Range("E42").Value = '  -1.81%  '

Set-TextValue "D43" '1.924'
$ws.Range("E43").Value = '  -2.64%  '

Set-TextValue "D44" '1.000'
$ws.Range("E44").Value = '  +0.02%  '

Set-TextValue "D45" '0.8362'

Set-TextValue "D46" '100.83'
$ws.Range("E46").Value = '  -0.02%  '

Set-TextValue "D47" '7.554'
$ws.Range("E47").Value = '  +0.59%  '

Set-TextValue "D48" '9.727'
$ws.Range("E48").Value = '  -0.67%  '

Set-TextValue "D49" '37.57'
$ws.Range("E49").Value = '  +0.81%  '

$ws.Range("D50").Value = '2.061.68'
$ws.Range("E50").Value = '  -1.03%  '

Set-TextValue "D51" '934.64'
$ws.Range("E51").Value = '  -5.13%  '
